# Apply the cryptos-list refresh described by the commit:
# "Updated cryptos list on Thu May 23 06:17:24 UTC 2024 with GitHub Actions"
#
# All of the edited cells hold plain TEXT in the source workbook (inline
# strings), even when the text looks numeric (e.g. "614.44", "0.0000253",
# "1.00"). Excel's COM Range.Value setter auto-coerces number-looking text
# into a real number, which would silently change the cell type. To keep
# the cells as text (matching the original authoring), each write:
#   1. forces the cell to Text format ("@") before assigning,
#   2. assigns the literal string value, then
#   3. restores the cell style to "Normal" so no stray number-format style
#      is left behind on the cell (keeps styles.xml/cell "s" attrs clean).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '69.423.68'
Set-TextValue 'E2' '  -0.75%  '
Set-TextValue 'D3' '3.762.74'
Set-TextValue 'E3' '  +0.06%  '
Set-TextValue 'E4' '  -0.01%  '
Set-TextValue 'D5' '614.44'
Set-TextValue 'E5' '  -1.09%  '
Set-TextValue 'D6' '176.98'
Set-TextValue 'E6' '  -2.47%  '
Set-TextValue 'D7' '3.760.76'
Set-TextValue 'E7' '  +0.03%  '
Set-TextValue 'E8' '  +0.07%  '
Set-TextValue 'D9' '0.526'
Set-TextValue 'E9' '  -1.48%  '
Set-TextValue 'E10' '  -1.80%  '
Set-TextValue 'D11' '6.42'
Set-TextValue 'E11' '  +1.94%  '
Set-TextValue 'E12' '  -1.55%  '
Set-TextValue 'D13' '39.64'
Set-TextValue 'E13' '  -4.01%  '
Set-TextValue 'D14' '0.0000253'
Set-TextValue 'E14' '  -2.30%  '
Set-TextValue 'D15' '4.390.74'
Set-TextValue 'E15' '  +0.02%  '
Set-TextValue 'D16' '3.761.78'
Set-TextValue 'E16' '  -0.18%  '
Set-TextValue 'D17' '69.514.00'
Set-TextValue 'E17' '  -0.78%  '
Set-TextValue 'D18' '7.52'
Set-TextValue 'E18' '  -1.03%  '
Set-TextValue 'E19' '  -3.53%  '
Set-TextValue 'D20' '508.03'
Set-TextValue 'D21' '16.49'
Set-TextValue 'E21' '  -1.77%  '
Set-TextValue 'D22' '9.49'
Set-TextValue 'E22' '  -0.81%  '
Set-TextValue 'D23' '0.731'
Set-TextValue 'E23' '  +0.67%  '
Set-TextValue 'D24' '2.46'
Set-TextValue 'E24' '  -2.08%  '
Set-TextValue 'D25' '86.22'
Set-TextValue 'E25' '  -1.01%  '
Set-TextValue 'B26' 'InternetComputer(DFINITY)'
Set-TextValue 'C26' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue 'D26' '12.80'
Set-TextValue 'E26' '  -2.76%  '
Set-TextValue 'B27' 'PEPE'
Set-TextValue 'C27' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D27' '0.0000140'
Set-TextValue 'E27' '  +4.47%  '
Set-TextValue 'D28' '10.48'
Set-TextValue 'E28' '  -5.21%  '
Set-TextValue 'E29' '  -0.19%  '
Set-TextValue 'B30' 'PancakeSwap'
Set-TextValue 'C30' 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D30' '3.00'
Set-TextValue 'E30' '  +2.91%  '
Set-TextValue 'B31' 'ImmutableX'
Set-TextValue 'C31' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue 'D31' '2.51'
Set-TextValue 'E31' '  -0.05%  '
Set-TextValue 'D32' '8.08'
Set-TextValue 'E32' '  +1.91%  '
Set-TextValue 'D33' '30.88'
Set-TextValue 'E33' '  -0.52%  '
Set-TextValue 'E34' '  -0.43%  '
Set-TextValue 'D35' '0.999'
Set-TextValue 'E35' '  -0.15%  '
Set-TextValue 'E36' '  -2.17%  '
Set-TextValue 'D37' '6.10'
Set-TextValue 'E37' '  -1.21%  '
Set-TextValue 'E38' '  +5.48%  '
Set-TextValue 'D39' '0.339'
Set-TextValue 'E39' '  +0.79%  '
Set-TextValue 'D40' '469.05'
Set-TextValue 'E40' '  +9.84%  '
Set-TextValue 'D41' '2.06'
Set-TextValue 'E41' '  -3.03%  '
Set-TextValue 'D42' '2.99'
Set-TextValue 'E42' '  +5.91%  '
Set-TextValue 'D43' '49.77'
Set-TextValue 'E43' '  -0.86%  '
Set-TextValue 'D44' '43.88'
Set-TextValue 'E44' '  -2.91%  '
Set-TextValue 'D45' '8.55'
Set-TextValue 'E45' '  -2.03%  '
Set-TextValue 'D46' '2.937.33'
Set-TextValue 'E46' '  -2.57%  '
Set-TextValue 'D47' '0.0360'
Set-TextValue 'E47' '  -1.17%  '
Set-TextValue 'D48' '27.34'
Set-TextValue 'E48' '  -0.66%  '
Set-TextValue 'B49' 'USDe'
Set-TextValue 'C49' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D49' '1.00'
Set-TextValue 'E49' '  +0.05%  '
Set-TextValue 'B50' 'Monero'
Set-TextValue 'C50' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue 'D50' '138.98'
Set-TextValue 'E50' '  +1.18%  '
Set-TextValue 'D51' '2.45'
Set-TextValue 'E51' '  -1.59%  '
